$wb = $excel.ActiveWorkbook

# Work on the "HomePage" sheet (3rd sheet / index 3)
$ws = $wb.Worksheets.Item("HomePage")
$ws.Activate()

# Add header + value for new column C
$ws.Range("C1").Value = "AddToBasket"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Leading apostrophe forces text-with-quote-prefix, matching style s="4"
# (quotePrefix) and shared-string cell content "2" from the target file.
$ws.Range("C2").Value = "'2"

# Size column C to fit its contents (closest achievable to the
# target bestFit width of 19.140625 characters)
$ws.Columns.Item(3).ColumnWidth = 18.3

# Update selection to C5
$ws.Range("C5").Select() | Out-Null
